# Generate Report for Handback
# Refresh the handoff/handback timestamps for the first tracked file
# (39a821cc-6cac-4628-8906-3cbb6bfb80e1.md) across the Overview, zh-cn and
# de-de report sheets.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$overview.Range("G2").Value = "2016-10-27 09:31:48"

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("H2").Value = "2016-10-27 09:31:35"
$zhcn.Range("K2").Value = "2016-10-27 09:32:25"

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("H2").Value = "2016-10-27 09:32:42"
$dede.Range("K2").Value = "2016-10-27 09:32:42"
